$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.810.03"
$ws.Range("E2").Value = "  +4.25%  "
$ws.Range("D3").Value = "1.874.81"
$ws.Range("E3").Value = "  +3.18%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'277.10"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.5283"
$ws.Range("E7").Value = "  +3.93%  "
$ws.Range("D8").Value = "'0.3410"
$ws.Range("E8").Value = "  -3.39%  "
$ws.Range("D9").Value = "'0.06949"
$ws.Range("E9").Value = "  +4.42%  "
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("D11").Value = "'0.8042"
$ws.Range("E11").Value = "  -2.51%  "
$ws.Range("D12").Value = "'0.07724"
$ws.Range("E12").Value = "  -1.89%  "
$ws.Range("D13").Value = "1.804.17"
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("D14").Value = "'5.186"
$ws.Range("E14").Value = "  +2.30%  "
$ws.Range("D15").Value = "'90.30"
$ws.Range("E16").Value = "  +3.34%  "
$ws.Range("D17").Value = "'0.9996"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "'0.000008046"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D20").Value = "26.850.62"
$ws.Range("D21").Value = "2.077.45"
$ws.Range("E21").Value = "  +1.38%  "
$ws.Range("D22").Value = "'4.748"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("D24").Value = "'6.177"
$ws.Range("E24").Value = "  +1.38%  "
$ws.Range("D25").Value = "'2.395"
$ws.Range("E25").Value = "  +8.26%  "
$ws.Range("D26").Value = "'146.47"
$ws.Range("E26").Value = "  +3.04%  "
$ws.Range("D27").Value = "'17.35"
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("D29").Value = "'113.62"
$ws.Range("E29").Value = "  +3.71%  "
$ws.Range("D30").Value = "'4.355"
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("D31").Value = "'4.305"
$ws.Range("E31").Value = "  +1.75%  "
$ws.Range("D32").Value = "'0.08910"
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("D33").Value = "'0.04917"
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("E34").Value = "  +3.35%  "
$ws.Range("D35").Value = "'0.7271"
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").Value = "'2.870"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("E37").Value = "  +4.76%  "
$ws.Range("D38").Value = "'2.341"
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("D40").Value = "'0.5131"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("D41").Value = "'0.9564"
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("D42").Value = "'116.50"
$ws.Range("E42").Value = "  +5.60%  "
$ws.Range("D43").Value = "'6.160"
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("D44").Value = "'8.126"
$ws.Range("E44").Value = "  +1.26%  "
$ws.Range("D45").Value = "'0.9999"
$ws.Range("D46").Value = "'0.4467"
$ws.Range("E46").Value = "  -1.83%  "
$ws.Range("D47").Value = "'0.1339"
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("D48").Value = "'9.338"
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("D49").Value = "'36.23"
$ws.Range("E49").Value = "  -0.89%  "
$ws.Range("D50").Value = "'0.05941"
$ws.Range("E50").Value = "  +1.69%  "
$ws.Range("D51").Value = "'1.492"
$ws.Range("E51").Value = "  -0.24%  "
